# "added cables to network"
# Rework the cable_lib sheet: rename two headers to match pandapower's
# cable-library column names, drop the two unused trailing columns, and
# add two new custom cable rows with their electrical parameters.

$wb = $excel.ActiveWorkbook

$cable = $wb.Worksheets.Item("cable_lib")
$buses = $wb.Worksheets.Item("Buses")

# --- header renames -------------------------------------------------
$cable.Range("B1").Value = "r_ohm_per_km"
$cable.Range("E1").Value = "max_i_ka"

# --- drop the now-unused q_mm2 / alpha columns -----------------------
$cable.Range("F1:G1").EntireColumn.Delete()

# --- new cable rows ---------------------------------------------------
$cable.Range("A2").Value = "CustomCable1"
$cable.Range("B2").Value = 0.15
$cable.Range("C2").Value = 0.08
$cable.Range("D2").Value = 300
$cable.Range("E2").Value = 0.35

$cable.Range("A3").Value = "CustomCable2"
$cable.Range("B3").Value = 0.05
$cable.Range("C3").Value = 0.12
$cable.Range("D3").Value = 250
$cable.Range("E3").Value = 0.5

# wrap the text in the new data rows, like the header row
$cable.Range("A2:D3").WrapText = $true

# --- restore on-screen selections ------------------------------------
[void]$buses.Activate()
[void]$buses.Range("B25").Select()

[void]$cable.Activate()
[void]$cable.Range("E4").Select()
